# Update weekly triaged issues
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> Component, Bugs, Features, Total
$data = @{
    3  = @("Popover", 6, 0, 6)
    4  = @("ComboBox", 4, 2, 6)
    5  = @("Dropdown", 4, 1, 5)
    7  = @("Tree", 5, 0, 5)
    8  = @("Drawer", 5, 0, 5)
    9  = @("TagPicker", 2, 2, 4)
    10 = @("Nav", 3, 0, 3)
    11 = @("Skeleton", 2, 1, 3)
    12 = @("Tooltip", 3, 0, 3)
    13 = @("Toolbar", 2, 1, 3)
    14 = @("Dialog", 2, 1, 3)
    15 = @("Virtualizer", 2, 0, 2)
    16 = @("Table", 2, 0, 2)
    17 = @("Portal", 2, 0, 2)
    18 = @("Popup", 1, 0, 1)
    20 = @("Avatar", 1, 0, 1)
    21 = @("FocusTrapZone", 1, 0, 1)
    22 = @("Card", 1, 0, 1)
    23 = @("Button", 0, 1, 1)
    24 = @("InfoLabel", 1, 0, 1)
    27 = @("Slider", 1, 1, 2)
    28 = @("Switch", 1, 0, 1)
    29 = @("Input", 0, 1, 1)
    30 = @("Checkbox", 1, 0, 1)
    31 = @("Badge", 1, 0, 1)
    32 = @("Toast", 1, 0, 1)
    33 = @("SpinButton", 1, 0, 1)
    34 = @("SplitButton", 1, 0, 1)
    35 = @("Calendar", 0, 0, 0)
    36 = @("Keytip", 0, 0, 0)
    37 = @("Pickers", 0, 0, 0)
    38 = @("List", 0, 0, 0)
    39 = @("Coachmark", 0, 0, 0)
    41 = @("Carousel", 0, 0, 0)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
}
